# Made the modal add new study sets
# Appends two new flashcard rows (simulating two new "study set" entries
# created through the add-new-study-set modal) to the bottom of the
# existing flashcards table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New study set #1
$ws.Range("A27").Value = "2023-10-08 02:25:54 9_4_1179775"
$ws.Range("B27").Value = "test test"
$ws.Range("C27").Value = "test"

# New study set #2
$ws.Range("A28").Value = "2023-10-08 02:34:14 9_4_3835979"
$ws.Range("B28").Value = "test test"
$ws.Range("C28").Value = "test"
